$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1) Update Version and Date values (rows 3 and 8 in the pre-edit layout).
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# 2) Insert a new "Jurisdiction" property row right after "Contact" (row 10),
#    pushing the existing rows 11-19 down to 12-20. We do this manually
#    (copy formatting, then copy values) instead of Rows.Insert() so every
#    shifted cell keeps the same style (s="2") the source row already had.
for ($r = 19; $r -ge 11; $r--) {
    $src = $ws.Range("A" + $r + ":B" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":B" + ($r + 1))

    # Copy cell formatting (keeps style index 2, matching every data row).
    $src.Copy()
    $dst.PasteSpecial(-4122)

    # Now copy the actual cell contents (preserves text/number/bool typing),
    # clearing the destination first so a blank source cell truly blanks
    # the destination instead of leaving stale content behind.
    $dst.ClearContents()
    $src.Copy()
    $dst.PasteSpecial(-4163)
}

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
